$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Remove the last two date columns (10/10, 17/10) from the schedule table
$t.Columns.Item(14).Delete()
$t.Columns.Item(13).Delete()

# Update NRP / Nama for the row that was "152015005" / "Ainan Taqarra Yusuf"
$t.Cell(2, 1).Range.Text = "152013001"
$t.Cell(2, 2).Range.Text = "Mochamad Angga Anggriawan"

# Update NRP / Nama for the row that was "152015001" / "Farhan Rafiqi"
$t.Cell(3, 1).Range.Text = "152013002"
$t.Cell(3, 2).Range.Text = "Gian Permana"
